$wb = $excel.ActiveWorkbook

# --- Step 1: the existing sheet ("Sheet1") holds the present-value scratchwork.
# Add a brand-new sheet right after it that will keep that original content under
# the name "Present value".
$ws1 = $wb.Worksheets.Item(1)
$ws2 = $wb.Worksheets.Add($null, $ws1)
$ws2.Name = "Present value"

# --- Step 2: populate "Present value" with the original Sheet1 scratchwork
# (principal/rate inputs, the discount-factor table, and the annuity sums).
$ws2.Range("A1").Value = 1
$ws2.Range("B1").Value = 0.1

$ws2.Range("B2").Value = "A_{n}"
$ws2.Range("C2").Value = "A_{n,due}"

$ws2.Range("A3").Value = 1
$ws2.Range("B3").Formula = "=`$A`$1*(1+`$B`$1)^-A3"
$ws2.Range("C3").Formula = "=`$A`$1*(1+`$B`$1)^(1-A3)"

$ws2.Range("A4").Formula = "=A3+1"
$ws2.Range("B4").Formula = "=`$A`$1*(1+`$B`$1)^-A4"
$ws2.Range("C4").Formula = "=`$A`$1*(1+`$B`$1)^(1-A4)"

$ws2.Range("A5").Formula = "=A4+1"
$ws2.Range("B5").Formula = "=`$A`$1*(1+`$B`$1)^-A5"
$ws2.Range("C5").Formula = "=`$A`$1*(1+`$B`$1)^(1-A5)"

$ws2.Range("B6").Formula = "=SUM(B3:B5)"
$ws2.Range("C6").Formula = "=SUM(C3:C5)"

$ws2.Range("B3:C6").NumberFormat = "[$$-409]#,##0.00;[RED]\-[$$-409]#,##0.00"

# --- Step 3: wipe the original sheet clean and turn it into the new
# "Future value" scratchwork (same annuity, but compounding forward to year 5
# instead of discounting back to year 0).
$ws1.Range("A1:C6").Clear()

$ws1.Range("A1").Value = 100
$ws1.Range("B1").Value = 0.1

$ws1.Range("A3").Value = 1
$ws1.Range("B3").Formula = "=`$A`$1*(1+`$B`$1)^(`$A`$7-A3)"

$ws1.Range("A4").Formula = "=A3+1"
$ws1.Range("B4").Formula = "=`$A`$1*(1+`$B`$1)^(`$A`$7-A4)"

$ws1.Range("A5").Formula = "=A4+1"
$ws1.Range("B5").Formula = "=`$A`$1*(1+`$B`$1)^(`$A`$7-A5)"

$ws1.Range("A6").Formula = "=A5+1"
$ws1.Range("B6").Formula = "=`$A`$1*(1+`$B`$1)^(`$A`$7-A6)"

$ws1.Range("A7").Formula = "=A6+1"
$ws1.Range("B7").Formula = "=`$A`$1*(1+`$B`$1)^(`$A`$7-A7)"

$ws1.Range("B8").Formula = "=SUM(B3:B7)"

$ws1.Range("B3:B8").NumberFormat = "[$$-409]#,##0.00;[RED]\-[$$-409]#,##0.00"

$ws1.Name = "Future value"

# --- Step 4: view tweaks recorded by the diff (zoom level + last selections
# on each sheet, with "Future value" left as the active tab).
[void]$ws2.Range("C13").Select()
$ws2.Application.ActiveWindow.Zoom = 400

[void]$ws1.Activate()
[void]$ws1.Range("C6").Select()
$ws1.Application.ActiveWindow.Zoom = 400
